# Chiffres COVID-19 Valais - data corrections + new daily rows
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Correct a handful of previously mis-entered daily new-case counts ---
$ws.Range("C23").Value = 37
$ws.Range("C26").Value = 38
$ws.Range("C222").Value = 16
$ws.Range("C253").Value = 541
$ws.Range("C262").Value = 302
$ws.Range("C264").Value = 120
$ws.Range("C334").Value = 70
$ws.Range("C338").Value = 125
$ws.Range("C370").Value = 90
$ws.Range("C487").Value = 9
$ws.Range("C489").Value = 6
$ws.Range("C492").Value = 9

# --- Fill in the newly-reported days (rows 493-496) ---
$ws.Range("C493").Value = 9
$ws.Range("E493").Value = 3
$ws.Range("F493").Value = 2
$ws.Range("G493").Value = 3
$ws.Range("L493").Value = 0
$ws.Range("M493").Value = 0

$ws.Range("C494").Value = 1
$ws.Range("E494").Value = 2
$ws.Range("F494").Value = 2
$ws.Range("G494").Value = 4
$ws.Range("L494").Value = 0
$ws.Range("M494").Value = 0

$ws.Range("C495").Value = 9
$ws.Range("E495").Value = 2
$ws.Range("F495").Value = 2
$ws.Range("G495").Value = 4
$ws.Range("L495").Value = 0
$ws.Range("M495").Value = 0

$ws.Range("C496").Value = 3
$ws.Range("E496").Value = 2
$ws.Range("F496").Value = 2
$ws.Range("G496").Value = 4
$ws.Range("L496").Value = 0
$ws.Range("M496").Value = 0

# --- Restore the view's active selection in the frozen bottom-right pane ---
$ws.Range("A2").Select()
